$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.236.66"
$ws.Range("E2").Value = "  +4.04%  "
$ws.Range("D3").Value = "3.485.65"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.68"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.33"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +7.07%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  +4.55%  "
$ws.Range("E11").Value = "  +4.81%  "
$ws.Range("D12").Value = "4.080.28"
$ws.Range("E12").Value = "  +3.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.76"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +7.26%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "3.496.87"
$ws.Range("E15").Value = "  +3.82%  "
$ws.Range("E16").Value = "  +4.19%  "
$ws.Range("D17").Value = "63.288.43"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.29"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +3.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.31"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +6.01%  "
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "392.84"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +2.98%  "
$ws.Range("E22").Value = "  +3.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.28"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +9.08%  "
$ws.Range("D26").Value = "3.630.38"
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.187"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.88"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +10.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  +5.69%  "
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("E32").Value = "  +6.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +4.01%  "
$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.62"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +29.82%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.17"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.34"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +8.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "171.66"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("E39").Value = "  +9.81%  "
$ws.Range("D40").Value = "3.524.63"
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  +4.10%  "
$ws.Range("E43").Value = "  +8.02%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.48"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.50"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +3.95%  "
$ws.Range("E46").Value = "  +10.81%  "
$ws.Range("D47").Value = "2.611.83"
$ws.Range("E47").Value = "  +6.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.78"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +7.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +17.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.76"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0271"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +5.12%  "
